$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2024-10-14 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-15 Tuesday", 2) | Out-Null

# Update all table cells (positional replace, by row/col,
# since some source values repeat but map to different targets)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "71+8=79"
$tbl.Cell(1,2).Range.Text = "92-20=72"
$tbl.Cell(1,3).Range.Text = "42-20=22"
$tbl.Cell(1,4).Range.Text = "84-49=35"
$tbl.Cell(1,5).Range.Text = "53+2=55"
$tbl.Cell(2,1).Range.Text = "56-21=35"
$tbl.Cell(2,2).Range.Text = "74-61=13"
$tbl.Cell(2,3).Range.Text = "49-8=41"
$tbl.Cell(2,4).Range.Text = "76-74=2"
$tbl.Cell(2,5).Range.Text = "2+18=20"
$tbl.Cell(3,1).Range.Text = "93-45=48"
$tbl.Cell(3,2).Range.Text = "2+29=31"
$tbl.Cell(3,3).Range.Text = "14+55=69"
$tbl.Cell(3,4).Range.Text = "89-23=66"
$tbl.Cell(3,5).Range.Text = "42+6=48"
$tbl.Cell(4,1).Range.Text = "63+7=70"
$tbl.Cell(4,2).Range.Text = "86-9=77"
$tbl.Cell(4,3).Range.Text = "25+28=53"
$tbl.Cell(4,4).Range.Text = "29-22=7"
$tbl.Cell(4,5).Range.Text = "76-62=14"
$tbl.Cell(5,1).Range.Text = "40+7=47"
$tbl.Cell(5,2).Range.Text = "19+34=53"
$tbl.Cell(5,3).Range.Text = "90-54=36"
$tbl.Cell(5,4).Range.Text = "66+14=80"
$tbl.Cell(5,5).Range.Text = "69-64=5"
$tbl.Cell(6,1).Range.Text = "39+29=68"
$tbl.Cell(6,2).Range.Text = "12-2=10"
$tbl.Cell(6,3).Range.Text = "39+39=78"
$tbl.Cell(6,4).Range.Text = "46+5=51"
$tbl.Cell(6,5).Range.Text = "86-5=81"
$tbl.Cell(7,1).Range.Text = "29-5=24"
$tbl.Cell(7,2).Range.Text = "77-37=40"
$tbl.Cell(7,3).Range.Text = "63+27=90"
$tbl.Cell(7,4).Range.Text = "78-61=17"
$tbl.Cell(7,5).Range.Text = "60-2=58"
$tbl.Cell(8,1).Range.Text = "12+1=13"
$tbl.Cell(8,2).Range.Text = "86-45=41"
$tbl.Cell(8,3).Range.Text = "54+36=90"
$tbl.Cell(8,4).Range.Text = "85-4=81"
$tbl.Cell(8,5).Range.Text = "97-27=70"
$tbl.Cell(9,1).Range.Text = "35+28=63"
$tbl.Cell(9,2).Range.Text = "6+4=10"
$tbl.Cell(9,3).Range.Text = "30+9=39"
$tbl.Cell(9,4).Range.Text = "82-44=38"
$tbl.Cell(9,5).Range.Text = "37+57=94"
$tbl.Cell(10,1).Range.Text = "73-31=42"
$tbl.Cell(10,2).Range.Text = "35+20=55"
$tbl.Cell(10,3).Range.Text = "39-38=1"
$tbl.Cell(10,4).Range.Text = "46-5=41"
$tbl.Cell(10,5).Range.Text = "36-4=32"
$tbl.Cell(11,1).Range.Text = "13+69=82"
$tbl.Cell(11,2).Range.Text = "64+4=68"
$tbl.Cell(11,3).Range.Text = "10-9=1"
$tbl.Cell(11,4).Range.Text = "3+87=90"
$tbl.Cell(11,5).Range.Text = "68-47=21"
$tbl.Cell(12,1).Range.Text = "48-25=23"
$tbl.Cell(12,2).Range.Text = "64-35=29"
$tbl.Cell(12,3).Range.Text = "50-7=43"
$tbl.Cell(12,4).Range.Text = "83-43=40"
$tbl.Cell(12,5).Range.Text = "14+47=61"
$tbl.Cell(13,1).Range.Text = "29-18=11"
$tbl.Cell(13,2).Range.Text = "43+28=71"
$tbl.Cell(13,3).Range.Text = "44+53=97"
$tbl.Cell(13,4).Range.Text = "85-42=43"
$tbl.Cell(13,5).Range.Text = "74+23=97"
$tbl.Cell(14,1).Range.Text = "74-42=32"
$tbl.Cell(14,2).Range.Text = "19+41=60"
$tbl.Cell(14,3).Range.Text = "53+33=86"
$tbl.Cell(14,4).Range.Text = "74+24=98"
$tbl.Cell(14,5).Range.Text = "3+86=89"
$tbl.Cell(15,1).Range.Text = "18+1=19"
$tbl.Cell(15,2).Range.Text = "49-36=13"
$tbl.Cell(15,3).Range.Text = "73+22=95"
$tbl.Cell(15,4).Range.Text = "92-30=62"
$tbl.Cell(15,5).Range.Text = "64-56=8"
$tbl.Cell(16,1).Range.Text = "90-14=76"
$tbl.Cell(16,2).Range.Text = "3+95=98"
$tbl.Cell(16,3).Range.Text = "84-31=53"
$tbl.Cell(16,4).Range.Text = "50-3=47"
$tbl.Cell(16,5).Range.Text = "89-79=10"
$tbl.Cell(17,1).Range.Text = "91-43=48"
$tbl.Cell(17,2).Range.Text = "52-0=52"
$tbl.Cell(17,3).Range.Text = "55-19=36"
$tbl.Cell(17,4).Range.Text = "22+9=31"
$tbl.Cell(17,5).Range.Text = "75-29=46"
$tbl.Cell(18,1).Range.Text = "7+31=38"
$tbl.Cell(18,2).Range.Text = "83-2=81"
$tbl.Cell(18,3).Range.Text = "59+13=72"
$tbl.Cell(18,4).Range.Text = "49-26=23"
$tbl.Cell(18,5).Range.Text = "63-50=13"
$tbl.Cell(19,1).Range.Text = "55+38=93"
$tbl.Cell(19,2).Range.Text = "64-45=19"
$tbl.Cell(19,3).Range.Text = "66-1=65"
$tbl.Cell(19,4).Range.Text = "8+89=97"
$tbl.Cell(19,5).Range.Text = "33+27=60"
$tbl.Cell(20,1).Range.Text = "24-11=13"
$tbl.Cell(20,2).Range.Text = "68+14=82"
$tbl.Cell(20,3).Range.Text = "40+41=81"
$tbl.Cell(20,4).Range.Text = "83-9=74"
$tbl.Cell(20,5).Range.Text = "59-1=58"
